# Payslip template employee record update
# (config/employee data refresh tied to the Orchestrator config-file fix)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Employee ID (row 8)
$ws.Range("C8").Value = 18

# Employee name & telephone (rows 9-11)
$ws.Range("C9").Value  = "Jerica"
$ws.Range("C10").Value = "Thacker"
$ws.Range("C11").Value = "(916) 984-8382"

# Hours worked / expenses / exchange-rate-pay figures
$ws.Range("G9").Value  = 23
$ws.Range("I9").Value  = 1840
$ws.Range("I11").Value = 1840
